$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Jason Holder"

# Force the numeric-looking columns (runs, balls, fours, sixes, sr) to Text
# format so values like "47" or "100.00" are stored as text, not numbers.
$ws.Range("E2:I8").NumberFormat = "@"

$data = New-Object 'object[,]' 8,13
$data[0,0] = 'matchNo'
$data[0,1] = 'teamName'
$data[0,2] = 'batterName'
$data[0,3] = 'states'
$data[0,4] = 'runs'
$data[0,5] = 'balls'
$data[0,6] = 'fours'
$data[0,7] = 'sixes'
$data[0,8] = 'sr'
$data[0,9] = 'opponentTeamName'
$data[0,10] = 'venue'
$data[0,11] = 'date'
$data[0,12] = 'result'
$data[1,0] = '37th'
$data[1,1] = 'Sunrisers Hyderabad'
$data[1,2] = 'Jason Holder'
$data[1,3] = ''
$data[1,4] = '47'
$data[1,5] = '29'
$data[1,6] = '0'
$data[1,7] = '5'
$data[1,8] = '162.06'
$data[1,9] = 'Punjab Kings'
$data[1,10] = 'Sharjah'
$data[1,11] = 'September 25'
$data[1,12] = 'Punjab Kings won by 5 runs'
$data[2,0] = '44th'
$data[2,1] = 'Sunrisers Hyderabad'
$data[2,2] = 'Jason Holder'
$data[2,3] = 'c Chahar b Thakur'
$data[2,4] = '5'
$data[2,5] = '5'
$data[2,6] = '0'
$data[2,7] = '0'
$data[2,8] = '100.00'
$data[2,9] = 'Chennai Super Kings'
$data[2,10] = 'Sharjah'
$data[2,11] = 'September 30'
$data[2,12] = 'Super Kings won by 6 wickets (with 2 balls remaining)'
$data[3,0] = '33rd'
$data[3,1] = 'Sunrisers Hyderabad'
$data[3,2] = 'Jason Holder'
$data[3,3] = 'c Shaw b Patel'
$data[3,4] = '10'
$data[3,5] = '9'
$data[3,6] = '0'
$data[3,7] = '1'
$data[3,8] = '111.11'
$data[3,9] = 'Delhi Capitals'
$data[3,10] = 'Dubai (DSC)'
$data[3,11] = 'September 22'
$data[3,12] = 'Capitals won by 8 wickets (with 13 balls remaining)'
$data[4,0] = '49th'
$data[4,1] = 'Sunrisers Hyderabad'
$data[4,2] = 'Jason Holder'
$data[4,3] = 'c Iyer b Varun'
$data[4,4] = '2'
$data[4,5] = '9'
$data[4,6] = '0'
$data[4,7] = '0'
$data[4,8] = '22.22'
$data[4,9] = 'Kolkata Knight Riders'
$data[4,10] = 'Dubai (DSC)'
$data[4,11] = 'October 03'
$data[4,12] = 'KKR won by 6 wickets (with 2 balls remaining)'
$data[5,0] = '52nd'
$data[5,1] = 'Sunrisers Hyderabad'
$data[5,2] = 'Jason Holder'
$data[5,3] = 'c Christian b Patel'
$data[5,4] = '16'
$data[5,5] = '13'
$data[5,6] = '2'
$data[5,7] = '0'
$data[5,8] = '123.07'
$data[5,9] = 'Royal Challengers Bangalore'
$data[5,10] = 'Abu Dhabi'
$data[5,11] = 'October 06'
$data[5,12] = 'Sunrisers won by 4 runs'
$data[6,0] = '55th'
$data[6,1] = 'Sunrisers Hyderabad'
$data[6,2] = 'Jason Holder'
$data[6,3] = 'c Boult b Coulter-Nile'
$data[6,4] = '1'
$data[6,5] = '2'
$data[6,6] = '0'
$data[6,7] = '0'
$data[6,8] = '50.00'
$data[6,9] = 'Mumbai Indians'
$data[6,10] = 'Abu Dhabi'
$data[6,11] = 'October 08'
$data[6,12] = 'Mumbai won by 42 runs'
$data[7,0] = '6th'
$data[7,1] = 'Sunrisers Hyderabad'
$data[7,2] = 'Jason Holder'
$data[7,3] = 'c Christian b Mohammed Siraj'
$data[7,4] = '4'
$data[7,5] = '5'
$data[7,6] = '0'
$data[7,7] = '0'
$data[7,8] = '80.00'
$data[7,9] = 'Royal Challengers Bangalore'
$data[7,10] = 'Chennai'
$data[7,11] = 'April 14'
$data[7,12] = 'RCB won by 6 runs'

$ws.Range("A1:M8").Value = $data
Write-Host "done"
